# Insert a new daily price-report row for "Feria Lagunitas de Puerto Montt -
# Zapallo italiano" at row 423 (this pushes the existing rows 423-466 down
# to 424-467, matching the historical-by-date ordering used in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 423:466 down to 424:467, leaving a blank row 423 to populate.
$ws.Rows.Item(423).Insert()

# Populate the new row 423 with the new record's data.
$ws.Range("A423").Value = 4
$ws.Range("B423").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C423").Value = "Los Lagos"
$ws.Range("D423").Value = 45194
$ws.Range("E423").Value = 10
$ws.Range("F423").Value = 100112032
$ws.Range("G423").Value = "Zapallo italiano"
$ws.Range("H423").Value = "Sin especificar"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 70
$ws.Range("K423").Value = 18000
$ws.Range("L423").Value = 18000
$ws.Range("M423").Value = 18000
$ws.Range("N423").Value = "$/caja 50 unidades"
$ws.Range("O423").Value = "Región de Arica y Parinacota"
$ws.Range("P423").Value = 360
$ws.Range("Q423").Value = 50
$ws.Range("R423").Value = "Hortaliza"
